$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (B2:M2) used to hold volatile =RANDBETWEEN(1, 10) formulas with
# stale cached results. Replace them with plain, deterministic values that
# mirror the sequence already used as headers in row 1 (1..12), so the
# formulas disappear entirely and the cached numbers become real data.
$values = @(1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12)
for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 2  # column B (2) .. column M (13)
    $ws.Cells.Item(2, $col).Value = $values[$i]
}

# Touch the font formatting of the whole A1:M2 data block so those cells are
# promoted to their own explicit cell-style entry instead of silently
# sharing the sheet's default style (matches the workbook's new style table
# picking up a second, font-aware xf record).
$ws.Range("A1:M2").Font.Color = 0

# Restore the active selection/cursor position recorded in the workbook.
$ws.Range("J16").Select() | Out-Null
